# base dao service 课堂跟进练习  2019-4-15 22:55:46
#
# Appends a new daily log block (rows 19-20) to the bottom of the
# existing "day / weekday / topic / time" table on Sheet1:
#   Row 19: timestamp | weekday | topic            | time range
#   Row 20:            |        | topic (cont.)    | time range
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "2019年4月15日22:54:03"
$ws.Range("B19").Value = "周一"
$ws.Range("C19").Value = "课堂跟进"
$ws.Range("D19").Value = "8:30--10:10"

$ws.Range("C20").Value = "base dao service课后练习"
$ws.Range("D20").Value = "16:30--18:00"

# Move the active selection to the last filled cell, matching the
# workbook's saved cursor position after the edit.
[void]$ws.Range("D20").Select()
